# Updated cryptos list with refreshed price/volume data (columns D and E).
#
# The new text values for several "Price" cells happen to look like plain
# numbers (e.g. "1.00", "322.01"). Excel's COM layer auto-detects those and
# silently coerces them to the Number type (losing the original text
# formatting and drifting float precision), so each cell is explicitly
# forced to Text format before the write and then restored to the sheet's
# normal (unstyled) look afterwards, matching how the source data is stored
# as plain inline/shared strings with no special cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    'D2' = '47.529.67'
    'E2' = '  +0.01%  '
    'D3' = '2.492.55'
    'D4' = '1.00'
    'E4' = '  +0.05%  '
    'D5' = '322.01'
    'E5' = '  -0.09%  '
    'D6' = '109.39'
    'E6' = '  +4.06%  '
    'E7' = '  -0.37%  '
    'E8' = '  -0.04%  '
    'E9' = '  +0.42%  '
    'D10' = '39.45'
    'E10' = '  +4.39%  '
    'E11' = '  -0.39%  '
    'E12' = '  +0.74%  '
    'D13' = '18.65'
    'E13' = '  +2.30%  '
    'E14' = '  +0.92%  '
    'D15' = '2.881.69'
    'E15' = '  +0.18%  '
    'D16' = '2.497.05'
    'E16' = '  -0.19%  '
    'E17' = '  +0.72%  '
    'D18' = '47.410.67'
    'E18' = '  +0.05%  '
    'E19' = '  +6.32%  '
    'E20' = '  +1.35%  '
    'E21' = '  +0.73%  '
    'D22' = '2.76'
    'E22' = '  +16.15%  '
    'E23' = '  +0.09%  '
    'D24' = '247.16'
    'E24' = '  -1.33%  '
    'E25' = '  -0.10%  '
    'E26' = '  -0.02%  '
    'D27' = '25.78'
    'E27' = '  -1.33%  '
    'D28' = '10.00'
    'E28' = '  +0.03%  '
    'D29' = '2.20'
    'E29' = '  -3.59%  '
    'E30' = '  +4.07%  '
    'D31' = '34.79'
    'E31' = '  -0.38%  '
    'D32' = '49.94'
    'E32' = '  +0.99%  '
    'D33' = '20.41'
    'E33' = '  +3.06%  '
    'E34' = '  -0.30%  '
    'E35' = '  +1.06%  '
    'D36' = '1.01'
    'E36' = '  +0.19%  '
    'E37' = '  +2.56%  '
    'E38' = '  +1.26%  '
    'E39' = '  -1.61%  '
    'E40' = '  +0.46%  '
    'D41' = '22.47'
    'E41' = '  +7.07%  '
    'E42' = '  -2.11%  '
    'D43' = '119.18'
    'E43' = '  -1.67%  '
    'E44' = '  +0.35%  '
    'D45' = '1.995.85'
    'E45' = '  +1.71%  '
    'E46' = '  +2.53%  '
    'E47' = '  -2.62%  '
    'D48' = '1.79'
    'E48' = '  -0.41%  '
    'D49' = '9.08'
    'E49' = '  -1.12%  '
    'D50' = '5.21'
    'E50' = '  -1.53%  '
    'E51' = '  +3.71%  '
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    $cell.Style = "Normal"
}

